$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# Update time_taken values in the "data" sheet (F2:F8)
$data.Range("F2").Value = "2021-10-05 14:33:49.439059"
$data.Range("F3").Value = "2021-10-05 14:33:49.439067"
$data.Range("F4").Value = "2021-10-05 14:33:49.439070"
$data.Range("F5").Value = "2021-10-05 14:33:49.439073"
$data.Range("F6").Value = "2021-10-05 14:33:49.439076"
$data.Range("F7").Value = "2021-10-05 14:33:49.439079"
$data.Range("F8").Value = "2021-10-05 14:33:49.439081"

# Add a new "metadata" worksheet right after "data"
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Match the page margins used on the "data" sheet (Excel's default-new-sheet
# margins differ from the ones already baked into this workbook).
$meta.PageSetup.LeftMargin = 54
$meta.PageSetup.RightMargin = 54
$meta.PageSetup.TopMargin = 72
$meta.PageSetup.BottomMargin = 72
$meta.PageSetup.HeaderMargin = 36
$meta.PageSetup.FooterMargin = 36

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Episodic Ataxia"
$meta.Range("C2").Value = 3179
$meta.Range("E2").Value = "2020-09-13T08:12:49.188761Z"
$meta.Range("F2").Value = "2021-10-05 14:33:49.435004"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3179/?format=json"

# D2 ("0.23") must land as literal text, not be auto-coerced to a number.
# Stage it as text in a scratch cell, copy just the value across, then wipe
# the scratch cell so it leaves no trace in the sheet.
$meta.Range("Z1").NumberFormat = "@"
$meta.Range("Z1").Value = "0.23"
$meta.Range("Z1").Copy()
$meta.Range("D2").PasteSpecial(-4123)
$meta.Range("Z1").Clear()
$excel.CutCopyMode = $false

# Copy the header/index cell formatting from the "data" sheet (bold, bordered,
# centered, top-aligned) onto the new sheet's header row and A2, matching the
# style already used for row 1 and column A on "data".
$data.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Keep the original "data" sheet as the active tab (matches the source doc)
$data.Activate()
